$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 378.33334
$ws.Range("I6").Value = 735
$ws.Range("K6").Value = 2205
$ws.Range("M6").Value = -2093
$ws.Range("H28").Value = 826.2857
$ws.Range("I28").Value = 796.3333
$ws.Range("J28").Value = 1006
$ws.Range("K28").Value = 796.3333
$ws.Range("L28").Value = 1006
$ws.Range("M28").Value = -311.3333
$ws.Range("N28").Value = -1976
$ws.Range("H45").Value = 998
$ws.Range("J45").Value = 998
$ws.Range("L45").Value = 2994
$ws.Range("N45").Value = -3378
$ws.Range("H74").Value = 4888.8887
$ws.Range("I74").Value = 3666.6667
$ws.Range("K74").Value = 3666.6667
$ws.Range("M74").Value = -2730.6667
$ws.Range("H77").Value = 4888.8887
$ws.Range("I77").Value = 3666.6667
$ws.Range("K77").Value = 18333.3335
$ws.Range("M77").Value = -13653.3335
$ws.Range("H80").Value = 17164.666
$ws.Range("J80").Value = 25499.5
$ws.Range("L80").Value = 76498.5
$ws.Range("N80").Value = -78494.5
$ws.Range("H83").Value = 17164.666
$ws.Range("J83").Value = 25499.5
$ws.Range("L83").Value = 229495.5
$ws.Range("N83").Value = -239479.5
$ws.Range("H135").Value = 2162.125
$ws.Range("I135").Value = 1899.5
$ws.Range("K135").Value = 17095.5
$ws.Range("M135").Value = -14560.5
$ws.Range("H137").Value = 4749.5
$ws.Range("I137").Value = 4500
$ws.Range("K137").Value = 13500
$ws.Range("M137").Value = -10950
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 15000
$ws.Range("I8").Value = 15000
$ws.Range("K8").Value = 15000
$ws.Range("M8").Value = -14856
$ws.Range("H45").Value = 3886
$ws.Range("I45").Value = 3651.2222
$ws.Range("J45").Value = 5999
$ws.Range("K45").Value = 3651.2222
$ws.Range("L45").Value = 5999
$ws.Range("M45").Value = -3274.2222
$ws.Range("N45").Value = -6753
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H74").Value = 3160.625
$ws.Range("I74").Value = 1547.5
$ws.Range("K74").Value = 1547.5
$ws.Range("M74").Value = -673.5
$ws.Range("H77").Value = 3160.625
$ws.Range("I77").Value = 1547.5
$ws.Range("K77").Value = 7737.5
$ws.Range("M77").Value = -3369.5
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H132").Value = 2337.6
$ws.Range("I132").Value = 1672
$ws.Range("K132").Value = 5016
$ws.Range("M132").Value = -2486
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 500
$ws.Range("I10").Value = 500
$ws.Range("K10").Value = 500
$ws.Range("M10").Value = -360
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H120").Value = 99988
$ws.Range("J120").Value = 99988
$ws.Range("L120").Value = 99988
$ws.Range("N120").Value = -109664
$ws.Range("H139").Value = 99995
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 99995
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 99995
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -110275
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 93.181816
$ws.Range("I7").Value = 71.77778000000001
$ws.Range("J7").Value = 189.5
$ws.Range("K7").Value = 71.77778000000001
$ws.Range("L7").Value = 189.5
$ws.Range("M7").Value = 41.22221999999999
$ws.Range("N7").Value = -415.5
$ws.Range("H22").Value = 299.5
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = 50
$ws.Range("H23").Value = 5000
$ws.Range("I23").Value = 5000
$ws.Range("K23").Value = 5000
$ws.Range("M23").Value = -4760
$ws.Range("H27").Value = 5000
$ws.Range("I27").Value = 5000
$ws.Range("K27").Value = 5000
$ws.Range("M27").Value = -4808
$ws.Range("H31").Value = 6293.6924
$ws.Range("J31").Value = 7076
$ws.Range("L31").Value = 7076
$ws.Range("N31").Value = -7666
$ws.Range("H34").Value = 6293.6924
$ws.Range("J34").Value = 7076
$ws.Range("L34").Value = 7076
$ws.Range("N34").Value = -7480
$ws.Range("H60").Value = 17018.092
$ws.Range("I60").Value = 13249.75
$ws.Range("J60").Value = 19171.428
$ws.Range("K60").Value = 13249.75
$ws.Range("L60").Value = 19171.428
$ws.Range("M60").Value = -12738.75
$ws.Range("N60").Value = -20193.428
$ws.Range("H95").Value = 18040.666
$ws.Range("J95").Value = 18040.666
$ws.Range("L95").Value = 18040.666
$ws.Range("N95").Value = -23532.666
$ws.Range("H96").Value = 3977
$ws.Range("J96").Value = 3977
$ws.Range("L96").Value = 3977
$ws.Range("N96").Value = -9469
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 699.6667
$ws.Range("I97").Value = 699.5
$ws.Range("J97").Value = 700
$ws.Range("K97").Value = 2098.5
$ws.Range("L97").Value = 2100
$ws.Range("M97").Value = -1602.5
$ws.Range("N97").Value = -3092
$ws.Range("H107").Value = 750
$ws.Range("J107").Value = 1000
$ws.Range("L107").Value = 3000
$ws.Range("N107").Value = -6840
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H11").Value = 12251125
$ws.Range("I11").Value = 14000000
$ws.Range("J11").Value = 9000
$ws.Range("K11").Value = 14000000
$ws.Range("L11").Value = 9000
$ws.Range("M11").Value = -13999861
$ws.Range("N11").Value = -9278
$ws.Range("H20").Value = 10009950
$ws.Range("I20").Value = 20000000
$ws.Range("J20").Value = 19900
$ws.Range("K20").Value = 20000000
$ws.Range("L20").Value = 19900
$ws.Range("M20").Value = -19999755
$ws.Range("N20").Value = -20390
$ws.Range("H24").Value = 20000000
$ws.Range("I24").Value = 20000000
$ws.Range("K24").Value = 20000000
$ws.Range("M24").Value = -19999827
$ws.Range("H102").Value = 6850.1113
$ws.Range("I102").Value = 6701.8667
$ws.Range("J102").Value = 7591.3335
$ws.Range("K102").Value = 6701.8667
$ws.Range("L102").Value = 7591.3335
$ws.Range("M102").Value = -5079.8667
$ws.Range("N102").Value = -10835.3335
$ws.Range("H113").Value = 1750
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 1255.5555
$ws.Range("J126").Value = 1287.5
$ws.Range("L126").Value = 3862.5
$ws.Range("N126").Value = -8802.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 1001199.75
$ws.Range("I12").Value = 1334599.6
$ws.Range("K12").Value = 1334599.6
$ws.Range("M12").Value = -1334429.6
$ws.Range("H16").Value = 1333.6666
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 1000.5
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 1000.5
$ws.Range("M16").Value = -1830
$ws.Range("N16").Value = -1340.5
$ws.Range("H20").Value = 16201.667
$ws.Range("J20").Value = 18800
$ws.Range("L20").Value = 18800
$ws.Range("N20").Value = -19252
$ws.Range("H22").Value = 1381.6666
$ws.Range("I22").Value = 1100
$ws.Range("J22").Value = 1945
$ws.Range("K22").Value = 1100
$ws.Range("L22").Value = 1945
$ws.Range("M22").Value = -805
$ws.Range("N22").Value = -2535
$ws.Range("H24").Value = 19900
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 19900
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 19900
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -20586
$ws.Range("H27").Value = 1381.6666
$ws.Range("I27").Value = 1100
$ws.Range("J27").Value = 1945
$ws.Range("K27").Value = 1100
$ws.Range("L27").Value = 1945
$ws.Range("M27").Value = -993
$ws.Range("N27").Value = -2159
$ws.Range("H74").Value = 20197
$ws.Range("I74").Value = 20197
$ws.Range("K74").Value = 20197
$ws.Range("M74").Value = -19199
$ws.Range("H77").Value = 20197
$ws.Range("I77").Value = 20197
$ws.Range("K77").Value = 60591
$ws.Range("M77").Value = -55599
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 19900
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 19900
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 19900
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -20476
$ws.Range("H18").Value = 18800
$ws.Range("J18").Value = 18800
$ws.Range("L18").Value = 18800
$ws.Range("N18").Value = -19146
$ws.Range("H20").Value = 21958.25
$ws.Range("J20").Value = 21958.25
$ws.Range("L20").Value = 21958.25
$ws.Range("N20").Value = -22438.25
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H46").Value = 47666.668
$ws.Range("J46").Value = 47666.668
$ws.Range("L46").Value = 47666.668
$ws.Range("N46").Value = -48128.668
$ws.Range("H134").Value = 47666.668
$ws.Range("J134").Value = 47666.668
$ws.Range("L134").Value = 143000.004
$ws.Range("N134").Value = -148070.004
